$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "About" sheet: drop the trailing "HK Notes" block (rows 10-12)
# ------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Rows("10:12").Delete()

# ------------------------------------------------------------------
# "Data" sheet: drop the HK-specific capacity/scale-factor block
# (rows 22-28: US total capacity / HK total capacity / scale factor)
# ------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")
$wsData.Rows("22:28").Delete()

# B17 reverts to the plain scientific-number style (no highlight fill),
# matching the formatting used by B10:B12 above it.
$wsData.Range("B10").Copy()
$wsData.Range("B17").PasteSpecial(-4122)   # xlPasteFormats

# ------------------------------------------------------------------
# "BTC" sheet: the 2010 starting capacity (B2) no longer gets scaled
# by the HK/US capacity ratio (Data!$B$28) - it now just references
# Data!B12 directly, matching the style used by the rest of row 2.
# ------------------------------------------------------------------
$wsBTC = $wb.Worksheets.Item("BTC")
$wsBTC.Range("C2").Copy()
$wsBTC.Range("B2").PasteSpecial(-4122)     # xlPasteFormats
$wsBTC.Range("B2").Formula = "=Data!B12"

$excel.CutCopyMode = 0
